$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2023"

$ws.Range("B11").Value = "Introduction to complex variable functions and their applications. Present differential equations of interest in physical engineering and develop solution techniques, verifying properties and resolution methods. Study of special functions in Physical Engineering."
$ws.Range("C11").Value = "Introduction to complex variable functions and their applications. Present differential equations of interest in physical engineering and develop solution techniques, verifying properties and resolution methods. Study of special functions in Physical Engineering."

$ws.Range("B13").Value = "Functions of a complex variable. Delta function. Partial differential equations in physical engineering: solution methods, solving boundary value problems, applications. Fourier Series and Integral Transforms. Special functions."
$ws.Range("C13").Value = "Functions of a complex variable. Delta function. Partial differential equations in physical engineering: solution methods, solving boundary value problems, applications. Fourier Series and Integral Transforms. Special functions."

$ws.Range("B15").Value = "Functions of a complex variable: infinite series, analytical functions, Cauchy Riemann conditions, boundary integrals, Cauchy's theorem, residue theorem, Delta function. Laplace equation, diffusion equation (of heat), wave equation (vibrating string), Fourier series, Fourier and Laplace integral transforms. Special functions: Legendre Polynomials, Spherical Harmonics, Bessel Functions."
$ws.Range("C15").Value = "Functions of a complex variable: infinite series, analytical functions, Cauchy Riemann conditions, boundary integrals, Cauchy's theorem, residue theorem, Delta function. Laplace equation, diffusion equation (of heat), wave equation (vibrating string), Fourier series, Fourier and Laplace integral transforms. Special functions: Legendre Polynomials, Spherical Harmonics, Bessel Functions."

$ws.Range("B19").Value = "Duas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + P2)/2"
$ws.Range("C19").Value = "Duas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + P2)/2"
